# fix a bug of order of sort
#
# Sheet "baseline" (sheet1) header row: simplify three rich-text header
# cells into plain text (two of them also get new wording to match the
# "treejoin" sheet's header naming), and give them the same
# "Droid Sans Fallback" font used elsewhere for the Chinese headers.
#
# Sheet "treejoin" (sheet2): the B column ("string edit distance" row
# counts before dedup) previously held a constant placeholder
# (1096209) for every row - replace it with the real per-row counts.
# Column E ("tree edit distance" flag) was always 1 - fix it to 0 for
# every row (bug in the original sort/merge order), which also ripples
# into the H total formula (E+F+G) recalculating automatically.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "baseline"
$ws2 = $wb.Worksheets.Item(2)   # "treejoin"

# ---------------------------------------------------------------------
# Sheet1 ("baseline") header row fixes
# ---------------------------------------------------------------------
$ws1.Range("B1").Value = "passjoin过滤后结果数"
$ws1.Range("B1").Font.Name = "Droid Sans Fallback"

$ws1.Range("C1").Value = "string edit distance过滤后结果数"
$ws1.Range("C1").Font.Name = "Droid Sans Fallback"

$ws1.Range("D1").Value = "最终结果数"
$ws1.Range("D1").Font.Name = "Droid Sans Fallback"
$ws1.Range("D1").Font.Color = 3947580

$ws1.Range("G1").Font.Name = "Droid Sans Fallback"

$ws1.Columns.Item(3).ColumnWidth = 27.651360544217667

# ---------------------------------------------------------------------
# Sheet2 ("treejoin") data fixes
# ---------------------------------------------------------------------
$bValues = @(1081,1281,2003,2471,2893,3731,4651,6511,9243,12359,17607,24277,35283,50641,69357,88315,115533,143209,173935,208665)

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $bValues[$i]
    $ws2.Cells.Item($row, 5).Value = 0
}

# ---------------------------------------------------------------------
# Selection / view tidy-up (cosmetic, matches the recorded diff)
# ---------------------------------------------------------------------
$ws1.Range("C1").Select()
$ws2.Range("B21").Select()
